{"js": "// Update the worksheet date and all the division problems to the new values.\n// Replacements are looked up and captured as Range references FIRST (all\n// searches against the original text), and only then are the edits applied.\n// This avoids any ambiguity created by a replacement value that happens to\n// equal another cell's original (not-yet-processed) value, e.g. \"76\u00f79=\" is\n// both a target (from \"32\u00f77=\") and a source (-> \"49\u00f75=\") in this sheet.\nconst replacements = [\n  [\"2024-11-01 Friday\", \"2024-11-02 Saturday\"],\n  [\"68\u00f72=\", \"98\u00f77=\"],\n  [\"32\u00f77=\", \"76\u00f79=\"],\n  [\"97\u00f76=\", \"30\u00f73=\"],\n  [\"95\u00f74=\", \"78\u00f72=\"],\n  [\"68\u00f74=\", \"71\u00f75=\"],\n  [\"63\u00f75=\", \"40\u00f73=\"],\n  [\"25\u00f73=\", \"35\u00f77=\"],\n  [\"95\u00f78=\", \"26\u00f78=\"],\n  [\"62\u00f75=\", \"83\u00f79=\"],\n  [\"84\u00f75=\", \"97\u00f78=\"],\n  [\"98\u00f72=\", \"67\u00f79=\"],\n  [\"18\u00f73=\", \"99\u00f79=\"],\n  [\"84\u00f78=\", \"25\u00f74=\"],\n  [\"88\u00f79=\", \"58\u00f76=\"],\n  [\"60\u00f79=\", \"41\u00f78=\"],\n  [\"76\u00f79=\", \"49\u00f75=\"],\n  [\"74\u00f76=\", \"94\u00f75=\"],\n  [\"85\u00f73=\", \"53\u00f72=\"],\n  [\"87\u00f72=\", \"74\u00f78=\"],\n  [\"11\u00f79=\", \"21\u00f78=\"],\n  [\"65\u00f79=\", \"51\u00f76=\"],\n  [\"42\u00f73=\", \"12\u00f76=\"],\n  [\"69\u00f75=\", \"53\u00f73=\"],\n  [\"75\u00f77=\", \"87\u00f78=\"],\n];\n\n// Pass 1: issue every search against the original document text and keep\n// the resulting (yet-unresolved) ranges around.\nconst searchResults = replacements.map(([oldText]) =>\n  context.document.body.search(oldText, { matchCase: true, matchWholeWord: false })\n);\nsearchResults.forEach((r) => r.load(\"items\"));\nawait context.sync();\n\n// Pass 2: now that every search has been resolved against the pristine\n// text, perform the text substitutions.\nfor (let i = 0; i < replacements.length; i++) {\n  const [, newText] = replacements[i];\n  const items = searchResults[i].items;\n  for (let j = 0; j < items.length; j++) {\n    items[j].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Update the worksheet date and all the division problems to the new values.\n#\n# NB: This runtime's `Find.Execute` searches/replaces across the *whole*\n# document regardless of which Range its Find object is scoped to (it\n# ignores the owning Range's boundaries once Replace is requested), so a\n# document-wide search-and-replace cannot be used safely here: \"76\u00f79=\" is\n# both a target value (written when \"32\u00f77=\" is updated) and, independently,\n# a pre-existing source value elsewhere in the sheet (-> \"49\u00f75=\"). Instead,\n# each cell/paragraph is edited directly through its own Range.Text, which\n# is purely offset-based and is unaffected by edits made to other Ranges.\n\n$d = $word.ActiveDocument\n\nfunction Set-RangeText($range, $newText) {\n    # Both paragraph Ranges and table-cell Ranges include a trailing mark\n    # character (paragraph mark \"\\r\", or end-of-cell marker) that must be\n    # excluded before assigning new text, otherwise it gets clobbered.\n    $range.MoveEnd(1, -1) | Out-Null\n    $range.Text = $newText\n}\n\n# Title / date paragraph.\nSet-RangeText $d.Paragraphs.Item(1).Range \"2024-11-02 Saturday\"\n\n# The 5x20 table: the five data rows (1, 5, 9, 13, 17) each hold five\n# division problems; the remaining rows are blank spacer rows.\n$table = $d.Tables.Item(1)\n\n$cellEdits = @(\n    @{ Row = 1;  Col = 1; Old = \"68\u00f72=\"; New = \"98\u00f77=\" },\n    @{ Row = 1;  Col = 2; Old = \"32\u00f77=\"; New = \"76\u00f79=\" },\n    @{ Row = 1;  Col = 4; Old = \"97\u00f76=\"; New = \"30\u00f73=\" },\n    @{ Row = 1;  Col = 5; Old = \"95\u00f74=\"; New = \"78\u00f72=\" },\n\n    @{ Row = 5;  Col = 1; Old = \"68\u00f74=\"; New = \"71\u00f75=\" },\n    @{ Row = 5;  Col = 2; Old = \"63\u00f75=\"; New = \"40\u00f73=\" },\n    @{ Row = 5;  Col = 3; Old = \"25\u00f73=\"; New = \"35\u00f77=\" },\n    @{ Row = 5;  Col = 4; Old = \"95\u00f78=\"; New = \"26\u00f78=\" },\n    @{ Row = 5;  Col = 5; Old = \"62\u00f75=\"; New = \"83\u00f79=\" },\n\n    @{ Row = 9;  Col = 1; Old = \"84\u00f75=\"; New = \"97\u00f78=\" },\n    @{ Row = 9;  Col = 2; Old = \"98\u00f72=\"; New = \"67\u00f79=\" },\n    @{ Row = 9;  Col = 3; Old = \"18\u00f73=\"; New = \"99\u00f79=\" },\n    @{ Row = 9;  Col = 4; Old = \"84\u00f78=\"; New = \"25\u00f74=\" },\n    @{ Row = 9;  Col = 5; Old = \"88\u00f79=\"; New = \"58\u00f76=\" },\n\n    @{ Row = 13; Col = 1; Old = \"60\u00f79=\"; New = \"41\u00f78=\" },\n    @{ Row = 13; Col = 2; Old = \"76\u00f79=\"; New = \"49\u00f75=\" },\n    @{ Row = 13; Col = 3; Old = \"74\u00f76=\"; New = \"94\u00f75=\" },\n    @{ Row = 13; Col = 4; Old = \"85\u00f73=\"; New = \"53\u00f72=\" },\n    @{ Row = 13; Col = 5; Old = \"87\u00f72=\"; New = \"74\u00f78=\" },\n\n    @{ Row = 17; Col = 1; Old = \"11\u00f79=\"; New = \"21\u00f78=\" },\n    @{ Row = 17; Col = 2; Old = \"65\u00f79=\"; New = \"51\u00f76=\" },\n    @{ Row = 17; Col = 3; Old = \"42\u00f73=\"; New = \"12\u00f76=\" },\n    @{ Row = 17; Col = 4; Old = \"69\u00f75=\"; New = \"53\u00f73=\" },\n    @{ Row = 17; Col = 5; Old = \"75\u00f77=\"; New = \"87\u00f78=\" }\n    # Row 1, Col 3 (\"32\u00f72=\") is intentionally left untouched.\n)\n\nforeach ($edit in $cellEdits) {\n    $cell = $table.Cell($edit.Row, $edit.Col)\n    $current = $cell.Range.Text.TrimEnd([char]7, [char]13)\n    if ($current -ne $edit.Old) {\n        Write-Output \"WARNING: cell ($($edit.Row),$($edit.Col)) expected '$($edit.Old)' but found '$current'\"\n    }\n    Set-RangeText $cell.Range $edit.New\n}\n\nWrite-Output \"Done.\"\n"}
